# Add Indian MF 1st Stab
# This script reproduces a weekly "MarketBeat" rank roll-forward: nine new
# date columns are inserted right after column A (the analyst-name column),
# pushing all the previously existing date columns further to the right,
# and the nine freshly inserted columns are populated with the new
# snapshot's header dates and "UN" (unchanged) markers. Finally one
# analyst (Citigroup, row 19) gets a brand-new rating-event note recorded
# under the newly inserted "Aug_25" column (C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert nine new blank columns before the existing date columns
#    (old B:V shifts to K:AE, keeping all values/styles/formats intact).
# ---------------------------------------------------------------------
$ws.Columns("B:J").Insert()

# ---------------------------------------------------------------------
# 2) Populate the new header row (row 1) with the nine new date labels.
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "Sep_08"
$ws.Range("C1").Value = "Aug_25"
$ws.Range("D1").Value = "Aug_04"
$ws.Range("E1").Value = "Jul_23"
$ws.Range("F1").Value = "Jul_17"
$ws.Range("G1").Value = "Jul_07"
$ws.Range("H1").Value = "Jun_30"
$ws.Range("I1").Value = "Jun_24"
$ws.Range("J1").Value = "Jun_16"

# ---------------------------------------------------------------------
# 3) Fill in "UN" (unchanged) for every analyst row across the nine new
#    columns, limited to each row's original populated extent (rows 1-29
#    went through old column V, rows 30-31 through old column P, and rows
#    32-33 through old column G -- all now nine columns further right).
# ---------------------------------------------------------------------
for ($r = 2; $r -le 29; $r++) {
    $ws.Range("B" + $r + ":J" + $r).Value = "UN"
}
for ($r = 30; $r -le 31; $r++) {
    $ws.Range("B" + $r + ":J" + $r).Value = "UN"
}
for ($r = 32; $r -le 33; $r++) {
    $ws.Range("B" + $r + ":J" + $r).Value = "UN"
}

# ---------------------------------------------------------------------
# 4) Record the new rating event for Citigroup (row 19) under the new
#    "Aug_25" column (C), replacing its "UN" placeholder, and give it the
#    same highlighted fill used for the other mid-table rating-change
#    notes (e.g. the existing Zacks/Credit Suisse Group event cells).
# ---------------------------------------------------------------------
$ws.Range("C19").Value = "8/16/2019,Set Price Target,Buy,$28.00"
$ws.Range("C19").Interior.ColorIndex = 38
